$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row 1 values according to the diff
$ws.Range("A1").Value = "Mesa"
$ws.Range("B1").Value = "Producto"
$ws.Range("C1").Value = "Cantidad"
$ws.Range("D1").Value = "Precio"
$ws.Range("E1").Value = "Fecha_Hora"
$ws.Range("F1").Value = "Total"
$ws.Range("G1").Value = "Estado"
$ws.Range("H1").Value = "Categoría"
$ws.Range("I1").Value = "Metodo_Pago"
$ws.Range("J1").Value = "Referencia"

# Delete the data rows (rows 2 through 8)
$ws.Range("A2:J8").EntireRow.Delete()
